$d = $word.ActiveDocument
$d.Content.Find.Execute("Gen 575 Winter 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Gen 575 Spring 2022", 2)
